$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Restore/update the value in C10 (Likelihood column, row "R30") from 18 to 1.
$ws.Range("C10").Value = 1

